$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 0
$ws.Range("A73").Value = 0
$ws.Range("A138").Value = 1
$ws.Range("A156").Value = 0
$ws.Range("A299").Value = 1
$ws.Range("A309").Value = 1
$ws.Range("A368").Value = 0
$ws.Range("A370").Value = 0
$ws.Range("A538").Value = 1
$ws.Range("A545").Value = 1
$ws.Range("A556").Value = 0
$ws.Range("A574").Value = 0
$ws.Range("A575:A576").Value = 1
$ws.Range("A581").Value = 1
$ws.Range("A599").Value = 1
$ws.Range("A621").Value = 0
$ws.Range("A666").Value = 1
$ws.Range("A672").Value = 0
$ws.Range("A691").Value = 1
$ws.Range("A698").Value = 1
$ws.Range("A706").Value = 1
$ws.Range("A716").Value = 1
$ws.Range("A753").Value = 0
$ws.Range("A788").Value = 0
$ws.Range("A830").Value = 1
$ws.Range("A847").Value = 0
$ws.Range("A874").Value = 0
$ws.Range("A877").Value = 1
$ws.Range("A914").Value = 1
$ws.Range("A916").Value = 0
$ws.Range("A918").Value = 0
$ws.Range("A920:A925").Value = 0
$ws.Range("A928").Value = 0
$ws.Range("A930").Value = 1
$ws.Range("A952:A953").Value = 0
$ws.Range("A958").Value = 0
$ws.Range("A963").Value = 1
$ws.Range("A966").Value = 1
$ws.Range("A970").Value = 0
$ws.Range("A972").Value = 0
$ws.Range("A976").Value = 0
$ws.Range("A981").Value = 0
$ws.Range("A984").Value = 0
$ws.Range("A988").Value = 1
$ws.Range("A991").Value = 0
$ws.Range("A994").Value = 0
$ws.Range("A997").Value = 0
$ws.Range("A1005").Value = 1
$ws.Range("A1017").Value = 0
$ws.Range("A1019").Value = 0
$ws.Range("A1022").Value = 1
$ws.Range("A1023").Value = 0
$ws.Range("A1030").Value = 0
$ws.Range("A1031").Value = 1
$ws.Range("A1034").Value = 1
$ws.Range("A1035").Value = 0
$ws.Range("A1043").Value = 1
$ws.Range("A1045").Value = 0
$ws.Range("A1048").Value = 1
$ws.Range("A1050").Value = 0
$ws.Range("A1057").Value = 1
$ws.Range("A1060").Value = 1
$ws.Range("A1065:A1066").Value = 0
$ws.Range("A1067").Value = 1
$ws.Range("A1070:A1071").Value = 0
$ws.Range("A1080").Value = 1
$ws.Range("A1081").Value = 0
$ws.Range("A1082:A1090").Value = 1
$ws.Range("A1092:A1097").Value = 1
$ws.Range("A1099:A1102").Value = 1
$ws.Range("A1152").Value = 1
$ws.Range("A1286").Value = 1
$ws.Range("A1374").Value = 0
$ws.Range("A1440").Value = 0
$ws.Range("A1469").Value = 0
$ws.Range("A1475").Value = 1
$ws.Range("A1485").Value = 1
$ws.Range("A1629").Value = 0
$ws.Range("A1770").Value = 0
$ws.Range("A1776").Value = 1
$ws.Range("A1783").Value = 1
$ws.Range("A1795").Value = 1
$ws.Range("A1800").Value = 1
